$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mon Sep 11 14:03:40 EDT 2023"
$ws.Range("B3").Value = "Mon Sep 11 14:03:54 EDT 2023"
$ws.Range("B4").Value = "Mon Sep 11 14:04:08 EDT 2023"
$ws.Range("B5").Value = "Mon Sep 11 14:04:21 EDT 2023"
